$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:N2").Value = "N/A"
$ws.Range("G5:N5").Value = "N/A"
$ws.Range("G9:N9").Value = "N/A"
$ws.Range("G11:N11").Value = "N/A"
$ws.Range("G12:N12").Value = "N/A"
$ws.Range("G13:N13").Value = "N/A"
$ws.Range("G16:N16").Value = "N/A"
$ws.Range("G17:N17").Value = "N/A"
$ws.Range("G19:N19").Value = "N/A"
$ws.Range("G20:N20").Value = "N/A"
$ws.Range("G22:N22").Value = "N/A"
$ws.Range("G27:N27").Value = "N/A"
$ws.Range("G30:N30").Value = "N/A"
$ws.Range("G31:N31").Value = "N/A"
$ws.Range("G33:N33").Value = "N/A"
$ws.Range("G35:N35").Value = "N/A"
$ws.Range("G37:N37").Value = "N/A"
$ws.Range("G41:N41").Value = "N/A"
$ws.Range("G42:N42").Value = "N/A"
$ws.Range("G47:N47").Value = "N/A"
$ws.Range("G53:N53").Value = "N/A"
$ws.Range("G60:N60").Value = "N/A"
$ws.Range("G61:N61").Value = "N/A"
$ws.Range("G63:N63").Value = "N/A"
$ws.Range("G66:N66").Value = "N/A"
$ws.Range("G69:N69").Value = "N/A"
$ws.Range("G73:N73").Value = "N/A"
$ws.Range("G76:N76").Value = "N/A"
$ws.Range("G78:N78").Value = "N/A"
$ws.Range("G82:N82").Value = "N/A"
$ws.Range("G86:N86").Value = "N/A"
$ws.Range("G87:N87").Value = "N/A"
$ws.Range("G92:N92").Value = "N/A"
$ws.Range("G95:N95").Value = "N/A"
$ws.Range("G96:N96").Value = "N/A"
$ws.Range("G101:N101").Value = "N/A"
$ws.Range("G106:N106").Value = "N/A"
$ws.Range("G108:N108").Value = "N/A"
$ws.Range("G113:N113").Value = "N/A"
$ws.Range("G117:N117").Value = "N/A"
$ws.Range("G123:N123").Value = "N/A"
$ws.Range("G127:N127").Value = "N/A"
$ws.Range("G133:N133").Value = "N/A"
$ws.Range("G135:N135").Value = "N/A"
$ws.Range("G136:N136").Value = "N/A"
$ws.Range("G141:N141").Value = "N/A"
$ws.Range("G142:N142").Value = "N/A"
$ws.Range("G143:N143").Value = "N/A"
$ws.Range("G146:N146").Value = "N/A"
$ws.Range("G149:N149").Value = "N/A"
$ws.Range("G151:N151").Value = "N/A"
$ws.Range("G155:N155").Value = "N/A"
$ws.Range("G156:N156").Value = "N/A"
$ws.Range("G160:N160").Value = "N/A"
$ws.Range("G166:N166").Value = "N/A"
$ws.Range("G169:N169").Value = "N/A"
$ws.Range("G174:N174").Value = "N/A"
$ws.Range("G179:N179").Value = "N/A"
$ws.Range("G180:N180").Value = "N/A"
$ws.Range("G182:N182").Value = "N/A"
$ws.Range("G185:N185").Value = "N/A"
$ws.Range("G188:N188").Value = "N/A"
$ws.Range("G189:N189").Value = "N/A"
$ws.Range("G190:N190").Value = "N/A"
